# The deck had five near-duplicate "Was there an increase in stock price
# volatility for Apple and Samsung..." slides (positions 16-20). The commit
# removes the one at position 18 (the one whose chart picture still carried
# the old placeholder alt-text "A diagram with a red line..." / image11.png),
# leaving the later duplicates (and every other slide) to shift up by one.
$p = $ppt.ActivePresentation
$p.Slides.Item(18).Delete()
